$wb = $excel.ActiveWorkbook

# --- Fix typo in "informacion" sheet: A10 referenced the wrong ferreteria name ---
$infoSheet = $wb.Worksheets.Item("informacion")
$infoSheet.Range("A10").Value = "FERRETERIA EL IMAN (DON BETO)"

# --- Restore missing data in "productos" sheet for rows 199-210 (FERRETERIA LOS ANGELES) ---
$prodSheet = $wb.Worksheets.Item("productos")

$prodSheet.Range("A198").Value = "FERRETERIA LOS ANGELES"

$prodSheet.Range("A199").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B199").Value = "Cemento"
$prodSheet.Range("C199").Value = "Cemento Pacasmayo Fortimax"
$prodSheet.Range("D199").Value = "Pacasmayo"

$prodSheet.Range("A200").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B200").Value = "Cemento"
$prodSheet.Range("C200").Value = "Cemento Mochica MS"
$prodSheet.Range("D200").Value = "Mochica"

$prodSheet.Range("A201").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B201").Value = "Morteros"
$prodSheet.Range("C201").Value = "Rapimix Asentado"
$prodSheet.Range("D201").Value = "Rapimix Pacasmayo"

$prodSheet.Range("A202").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B202").Value = "Morteros"
$prodSheet.Range("C202").Value = "Rapimix Tarrajeo"
$prodSheet.Range("D202").Value = "Rapimix Pacasmayo"

$prodSheet.Range("A203").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B203").Value = "Morteros"
$prodSheet.Range("C203").Value = "Rapimix Concreto Seco"
$prodSheet.Range("D203").Value = "Rapimix Pacasmayo"

$prodSheet.Range("A204").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B204").Value = "Fierro"
$prodSheet.Range("C204").Value = "Fierro barra 5/8"
$prodSheet.Range("D204").Value = "Sider"

$prodSheet.Range("A205").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B205").Value = "Tubería"
$prodSheet.Range("C205").Value = "Tubo de agua 3/4"
$prodSheet.Range("D205").Value = "Tuboplast"

$prodSheet.Range("A206").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B206").Value = "Ladrillos"
$prodSheet.Range("C206").Value = "Ladrillo techo 12"
$prodSheet.Range("D206").Value = "El Roble"

$prodSheet.Range("A207").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B207").Value = "Ladrillos"
$prodSheet.Range("C207").Value = "Ladrillo techo 15"
$prodSheet.Range("D207").Value = "El Roble"

$prodSheet.Range("A208").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B208").Value = "Ladrillos"
$prodSheet.Range("C208").Value = "Ladrillo king kong"
$prodSheet.Range("D208").Value = "El Roble"

$prodSheet.Range("A209").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B209").Value = "Fierro"
$prodSheet.Range("C209").Value = "Fierro barra 1/2"
$prodSheet.Range("D209").Value = "Sider"

$prodSheet.Range("A210").Value = "FERRETERIA LOS ANGELES"
$prodSheet.Range("B210").Value = "Fierro"
$prodSheet.Range("C210").Value = "Fierro barra 6mm"
$prodSheet.Range("D210").Value = "Sider"
